$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update quantity for the remaining product row
$ws.Range("C2").Value = 32

# Remove the second product row entirely (DS2310WMUS-LF / 243)
$ws.Range("B3:C3").Delete()
